$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates -------------------------------------------

# "Datos actualizados..." timestamp banner (A1) bumped from 00:23 to 01:40.
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 01:40"

# Countries re-sorted by total cases (column B): a handful of neighbouring
# pairs swapped rank, so the country label in those two rows swaps while
# the underlying per-row stats are simply the new day's numbers below.
# Panama (row 39) overtook Oman (row 38):
$ws.Range("A38").Value = "Panama"
$ws.Range("A39").Value = "Oman"

# Surinam (row 121) overtook Somalia (row 120):
$ws.Range("A120").Value = "Surinam"
$ws.Range("A121").Value = "Somalia"

# Islas Malvinas (row 214) overtook Montserrat (row 213):
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Updated per-country statistics (columns B..H) -------------------------

# Row 4
$ws.Range("B4").Value = 5696643
$ws.Range("C4").Value = 40669
$ws.Range("D4").Value = 3054701
$ws.Range("E4").Value = 2465739
$ws.Range("G4").Value = 1129
$ws.Range("H4").Value = 176203

# Row 5
$ws.Range("B5").Value = 3460413
$ws.Range("C5").Value = 48541
$ws.Range("E5").Value = 733970
$ws.Range("G5").Value = 1170
$ws.Range("H5").Value = 111189

# Row 11
$ws.Range("B11").Value = 502178
$ws.Range("C11").Value = 13056
$ws.Range("D11").Value = 326298
$ws.Range("E11").Value = 159901
$ws.Range("G11").Value = 360
$ws.Range("H11").Value = 15979

# Row 16
$ws.Range("B16").Value = 312659
$ws.Range("C16").Value = 6693
$ws.Range("E16").Value = 77604
$ws.Range("G16").Value = 282
$ws.Range("H16").Value = 6330

# Row 22
$ws.Range("B22").Value = 229700
$ws.Range("C22").Value = 1595
$ws.Range("E22").Value = 16486

# Row 38
$ws.Range("B38").Value = 83754
$ws.Range("C38").Value = 964
$ws.Range("D38").Value = 58274
$ws.Range("E38").Value = 23653
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 1827

# Row 39
$ws.Range("B39").Value = 83606
$ws.Range("C39").Value = 188
$ws.Range("D39").Value = 78188
$ws.Range("E39").Value = 4815
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 603

# Row 48
$ws.Range("B48").Value = 57550
$ws.Range("C48").Value = 865
$ws.Range("D48").Value = 43835
$ws.Range("E48").Value = 12587
$ws.Range("G48").Value = 13
$ws.Range("H48").Value = 1128

# Row 53
$ws.Range("B53").Value = 47950
$ws.Range("C53").Value = 369
$ws.Range("D53").Value = 44278
$ws.Range("E53").Value = 3494

# Row 56
$ws.Range("B56").Value = 42325
$ws.Range("C56").Value = 179
$ws.Range("D56").Value = 35197
$ws.Range("E56").Value = 5630

# Row 74
$ws.Range("B74").Value = 20798
$ws.Range("C74").Value = 315
$ws.Range("E74").Value = 4779

# Row 81
$ws.Range("B81").Value = 14820
$ws.Range("C81").Value = 151
$ws.Range("D81").Value = 9931
$ws.Range("E81").Value = 4362
$ws.Range("G81").Value = 8
$ws.Range("H81").Value = 527

# Row 89
$ws.Range("B89").Value = 10162
$ws.Range("C89").Value = 51
$ws.Range("E89").Value = 1043

# Row 120
$ws.Range("B120").Value = 3295
$ws.Range("C120").Value = 79
$ws.Range("D120").Value = 2227
$ws.Range("E120").Value = 1014
$ws.Range("H120").Value = 54

# Row 121
$ws.Range("B121").Value = 3265
$ws.Range("C121").Value = 8
$ws.Range("D121").Value = 2396
$ws.Range("E121").Value = 776
$ws.Range("H121").Value = 93

# Row 129
$ws.Range("D129").Value = 1290
$ws.Range("E129").Value = 1157

# Row 143
$ws.Range("B143").Value = 1493
$ws.Range("C143").Value = 8
$ws.Range("D143").Value = 1228
$ws.Range("E143").Value = 225

# Row 158
$ws.Range("B158").Value = 994
$ws.Range("C158").Value = 5
$ws.Range("E158").Value = 436

# Row 166
$ws.Range("B166").Value = 686
$ws.Range("C166").Value = 57
$ws.Range("E166").Value = 534

# Row 213
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

Write-Output "Applied country/provincia updates"
